# Horarios Línea 141 - actualización a 18:10:41 (102 filas nuevas/movidas)
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Hoja 1: LP1912  (dimension A1:E77 -> A1:E84, Total filas 72 -> 79)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 18:10:41"
$ws1.Range("A3").Value = "Total filas: 79"

$rows1 = @(
    @(44, '16:52:37', '18:03', '17_ROMERO', 71, 'LP1912'),
    @(45, '17:35:09', '18:03', '23_HERNANDEZ', 28, 'LP1912'),
    @(46, '16:46:42', '18:04', '14_ABASTO', 78, 'LP1912'),
    @(47, '17:35:09', '18:05', '14_ABASTO', 30, 'LP1912'),
    @(48, '18:10:41', '18:11', '10_OLMOS', 1, 'LP1912'),
    @(49, '18:10:41', '18:11', '16_SANTA ANA', 1, 'LP1912'),
    @(50, '16:52:37', '18:14', '10_OLMOS', 82, 'LP1912'),
    @(51, '17:47:22', '18:21', '16_SANTA ANA', 34, 'LP1912'),
    @(52, '16:46:42', '18:24', '11_ETCHEVERRY', 98, 'LP1912'),
    @(53, '17:35:09', '18:25', '11_ETCHEVERRY', 50, 'LP1912'),
    @(54, '17:13:30', '18:27', '15_ABASTO', 74, 'LP1912'),
    @(55, '17:47:22', '18:31', '16_SANTA ANA', 44, 'LP1912'),
    @(56, '17:35:09', '18:31', '23_HERNANDEZ', 56, 'LP1912'),
    @(57, '17:54:43', '18:33', '23_HERNANDEZ', 39, 'LP1912'),
    @(58, '16:46:42', '18:34', '14X44_ABASTO', 108, 'LP1912'),
    @(59, '16:46:42', '18:38', '17X38_ROMERO', 112, 'LP1912'),
    @(60, '17:13:30', '18:41', '14_ABASTO', 88, 'LP1912'),
    @(61, '16:46:42', '18:41', '16_P MOR-SANTA ANA', 115, 'LP1912'),
    @(62, '17:47:22', '18:44', '14_ABASTO', 57, 'LP1912'),
    @(63, '17:35:09', '18:45', '14_ABASTO', 70, 'LP1912'),
    @(64, '17:35:09', '18:51', '15_ABASTO', 76, 'LP1912'),
    @(65, '17:54:43', '18:53', '16_SANTA ANA', 59, 'LP1912'),
    @(66, '17:35:09', '18:59', '10_OLMOS', 84, 'LP1912'),
    @(67, '17:13:30', '19:01', '17_ROMERO', 108, 'LP1912'),
    @(68, '18:10:41', '19:03', '23_HERNANDEZ', 53, 'LP1912'),
    @(69, '17:13:30', '19:11', '81_EL PELIGRO', 118, 'LP1912'),
    @(70, '18:10:41', '19:14', '14_ABASTO', 64, 'LP1912'),
    @(71, '17:47:22', '19:17', '27_EL RETIRO', 90, 'LP1912'),
    @(72, '17:35:09', '19:19', '27_EL RETIRO', 104, 'LP1912'),
    @(73, '17:54:43', '19:20', '215C_EL PATO', 86, 'LP1912'),
    @(74, '17:35:09', '19:21', '215C_EL PATO', 106, 'LP1912'),
    @(75, '17:35:09', '19:29', '225_GOMEZ', 114, 'LP1912'),
    @(76, '17:54:43', '19:30', '215_EL PELIGRO', 96, 'LP1912'),
    @(77, '17:54:43', '19:30', '27_EL RETIRO', 96, 'LP1912'),
    @(78, '17:35:09', '19:31', '215_EL PELIGRO', 116, 'LP1912'),
    @(79, '18:10:41', '19:31', '27_EL RETIRO', 81, 'LP1912'),
    @(80, '17:47:22', '19:40', '17X38_ROMERO', 113, 'LP1912'),
    @(81, '17:47:22', '19:44', '11_ETCHEVERRY', 117, 'LP1912'),
    @(82, '17:54:43', '19:51', '81_EL PELIGRO', 117, 'LP1912'),
    @(83, '18:10:41', '19:58', '14X44_ABASTO', 108, 'LP1912'),
    @(84, '18:10:41', '20:00', '215C_EL PATO', 110, 'LP1912')
)

foreach ($row in $rows1) {
    $r = $row[0]
    $ws1.Cells.Item($r, 1).Value = $row[1]
    $ws1.Cells.Item($r, 2).Value = $row[2]
    $ws1.Cells.Item($r, 3).Value = $row[3]
    $ws1.Cells.Item($r, 4).Value = $row[4]
    $ws1.Cells.Item($r, 5).Value = $row[5]
}

# ---------------------------------------------------------------------------
# Hoja 2: LP1912-215  (dimension A1:E15 -> A1:E16, Total filas 10 -> 11)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 18:10:41"
$ws2.Range("A3").Value = "Total filas: 11"

$ws2.Cells.Item(16, 1).Value = "18:10:41"
$ws2.Cells.Item(16, 2).Value = "20:00"
$ws2.Cells.Item(16, 3).Value = "215C_EL PATO"
$ws2.Cells.Item(16, 4).Value = 110
$ws2.Cells.Item(16, 5).Value = "LP1912"

# ---------------------------------------------------------------------------
# Hoja 3: 6203-6173  (dimension A1:E12 -> A1:E13, Total filas 7 -> 8)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 18:10:41"
$ws3.Range("A3").Value = "Total filas: 8"

$ws3.Cells.Item(13, 1).Value = "18:10:41"
$ws3.Cells.Item(13, 2).Value = "19:15"
$ws3.Cells.Item(13, 3).Value = "215B_LP-P MOR-1 Y 57"
$ws3.Cells.Item(13, 4).Value = 65
$ws3.Cells.Item(13, 5).Value = "L6173"
